# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the source file
# has moved from "In Translation" to "Ready for handoff", refreshes the
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps,
# and widens the Status columns to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# Column width (in Excel "characters" units) that renders to the same
# OOXML <col width> the report generator produced for the widened
# Status columns.
$statusColWidth = 16.333333333333332

# ---------------------------------------------------------------------
# Overview sheet: per-language status + latest generation timestamp
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-27 12:55:26"

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# ---------------------------------------------------------------------
# zh-cn sheet: Latest Handoff Datetime + widened Status column
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("H2").Value = "2016-08-27 12:55:22"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# ---------------------------------------------------------------------
# de-de sheet: Latest Handoff Datetime + widened Status column
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("H2").Value = "2016-08-27 12:55:26"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
